# Generate Report for Handoff
# Update status / handoff-datetime for the "ed31471e-b328-441f-a105-ede1d361a2df" file
# from "Handed back: in sync with en-US" to "Ready for handoff" across the Overview,
# zh-cn and de-de sheets, and bump the corresponding "Latest Handoff Datetime" values.

$wb = $excel.ActiveWorkbook

# --- Overview sheet ---
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("B3").Value = "Ready for handoff"
$wsOverview.Range("C3").Value = "Ready for handoff"
$wsOverview.Range("D3").Value = "2016-03-21 20:47:08"

# --- zh-cn sheet ---
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("C3").Value = "Ready for handoff"
$wsZhCn.Range("E3").Value = "2016-03-21 20:47:04"

# --- de-de sheet ---
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("C3").Value = "Ready for handoff"
$wsDeDe.Range("E3").Value = "2016-03-21 20:47:08"
